{"js": "// 1) Rewrite the \"For 4 years, he served...\" sentence into the new\n//    \"Bradley led the internal AWS community of practice at CapTech...\"\n//    sentence (also drops the trailing space that used to follow the\n//    final period).\nconst oldSentence =\n  \"For 4 years, he served as the internal lead for CapTech's AWS community \" +\n  \"of practice (+150 members) where he spread cloud expertise and created \" +\n  \"opportunities for others to grow. \";\nconst newSentence =\n  \"Bradley led the internal AWS community of practice at CapTech \" +\n  \"(+150 members) where he spread cloud expertise and created \" +\n  \"opportunities for others to grow.\";\n\nconst sentenceResults = context.document.body.search(oldSentence, { matchCase: true });\nsentenceResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < sentenceResults.items.length; i++) {\n  sentenceResults.items[i].insertText(newSentence, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Fix the \"Hashicorp\" spelling to \"HashiCorp\".\nconst hashicorpResults = context.document.body.search(\"Hashicorp\", { matchCase: true });\nhashicorpResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < hashicorpResults.items.length; i++) {\n  hashicorpResults.items[i].insertText(\"HashiCorp\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) Fix every \"Devops\" spelling to \"DevOps\".\nconst devopsResults = context.document.body.search(\"Devops\", { matchCase: true });\ndevopsResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < devopsResults.items.length; i++) {\n  devopsResults.items[i].insertText(\"DevOps\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Rewrite the \"For 4 years, he served...\" sentence into the new\n#    \"Bradley led the internal AWS community of practice at CapTech...\"\n#    sentence (also drops the trailing space that used to follow the\n#    final period).\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"For 4 years, he served as the internal lead for CapTech's AWS community of practice (+150 members) where he spread cloud expertise and created opportunities for others to grow. \"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"Bradley led the internal AWS community of practice at CapTech (+150 members) where he spread cloud expertise and created opportunities for others to grow.\"\n$find1.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 2) Fix the \"Hashicorp\" spelling to \"HashiCorp\".\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Hashicorp\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"HashiCorp\"\n$find2.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 3) Fix every \"Devops\" spelling to \"DevOps\".\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Text = \"Devops\"\n$find3.Replacement.ClearFormatting()\n$find3.Replacement.Text = \"DevOps\"\n$find3.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n"}
